$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format before writing, so decimal-looking
# strings (e.g. "551.46") are stored as text rather than numbers,
# matching the original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.052.79'
$ws.Range('E2').Value = '  -0.41%  '

$ws.Range('D3').Value = '3.513.95'
$ws.Range('E3').Value = '  -2.06%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = '203.12'
$ws.Range('E5').Value = '  +6.23%  '

$ws.Range('D6').Value = '551.46'
$ws.Range('E6').Value = '  -4.52%  '

$ws.Range('D7').Value = '3.505.47'
$ws.Range('E7').Value = '  -2.13%  '

$ws.Range('D8').Value = '0.600'
$ws.Range('E8').Value = '  -2.87%  '

$ws.Range('E9').Value = '  -0.20%  '

$ws.Range('D10').Value = '0.655'
$ws.Range('E10').Value = '  -3.38%  '

$ws.Range('D11').Value = '61.20'
$ws.Range('E11').Value = '  +12.03%  '

$ws.Range('D12').Value = '0.143'
$ws.Range('E12').Value = '  -5.06%  '

$ws.Range('D13').Value = '0.0000271'
$ws.Range('E13').Value = '  -1.15%  '

$ws.Range('D14').Value = '9.87'
$ws.Range('E14').Value = '  -0.52%  '

$ws.Range('D15').Value = '4.105.82'
$ws.Range('E15').Value = '  -1.29%  '

$ws.Range('D16').Value = '3.531.33'
$ws.Range('E16').Value = '  -1.56%  '

$ws.Range('D17').Value = '0.124'
$ws.Range('E17').Value = '  -0.91%  '

$ws.Range('D18').Value = '18.60'
$ws.Range('E18').Value = '  +1.35%  '

$ws.Range('D19').Value = '66.929.07'
$ws.Range('E19').Value = '  -0.40%  '

$ws.Range('D20').Value = '11.85'
$ws.Range('E20').Value = '  -3.86%  '

$ws.Range('E21').Value = '  -2.77%  '

$ws.Range('D22').Value = '389.49'
$ws.Range('E22').Value = '  -3.29%  '

$ws.Range('D23').Value = '4.01'
$ws.Range('E23').Value = '  -4.69%  '

$ws.Range('D24').Value = '11.88'
$ws.Range('E24').Value = '  -9.35%  '

$ws.Range('D25').Value = '82.64'
$ws.Range('E25').Value = '  -3.61%  '

$ws.Range('D26').Value = '6.13'
$ws.Range('E26').Value = '  +0.47%  '

$ws.Range('D27').Value = '2.80'
$ws.Range('E27').Value = '  -4.56%  '

$ws.Range('D28').Value = '12.01'
$ws.Range('E28').Value = '  -4.31%  '

$ws.Range('D29').Value = '3.73'
$ws.Range('E29').Value = '  -1.72%  '

$ws.Range('D30').Value = '8.88'
$ws.Range('E30').Value = '  -2.59%  '

$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '30.70'
$ws.Range('E31').Value = '  -1.69%  '

$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = '692.47'
$ws.Range('E32').Value = '  +3.79%  '

$ws.Range('D33').Value = '7.30'
$ws.Range('E33').Value = '  -9.70%  '

$ws.Range('D34').Value = '11.75'
$ws.Range('E34').Value = '  -3.66%  '

$ws.Range('D35').Value = '63.07'
$ws.Range('E35').Value = '  -1.57%  '

$ws.Range('D36').Value = '0.110'
$ws.Range('E36').Value = '  -4.34%  '

$ws.Range('D37').Value = '39.94'
$ws.Range('E37').Value = '  -6.44%  '

$ws.Range('D38').Value = '0.406'
$ws.Range('E38').Value = '  -3.78%  '

$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.05%  '

$ws.Range('D40').Value = '3.08'
$ws.Range('E40').Value = '  -1.54%  '

$ws.Range('D41').Value = '3.121.99'
$ws.Range('E41').Value = '  -0.66%  '

$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.13%  '

$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = '0.129'
$ws.Range('E43').Value = '  -3.57%  '

$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value = '0.0₃0708'
$ws.Range('E44').Value = '  -10.38%  '

$ws.Range('D45').Value = '2.85'
$ws.Range('E45').Value = '  +18.80%  '

$ws.Range('D46').Value = '2.51'
$ws.Range('E46').Value = '  -12.51%  '

$ws.Range('D47').Value = '2.72'
$ws.Range('E47').Value = '  +7.75%  '

$ws.Range('D48').Value = '0.0399'
$ws.Range('E48').Value = '  -4.62%  '

$ws.Range('E49').Value = '  -2.74%  '

$ws.Range('D50').Value = '2.99'
$ws.Range('E50').Value = '  -3.62%  '

$ws.Range('D51').Value = '136.95'
$ws.Range('E51').Value = '  -4.55%  '

# Restore default styling on column D (clears the "@" number format
# marker so untouched cells keep their original unstyled look).
$ws.Range("D2:D51").Style = "Normal"
